$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1030.579
$ws.Cells.Item(32, 9).Value = 865.3333
$ws.Cells.Item(32, 10).Value = 1106.8462
$ws.Cells.Item(32, 11).Value = 865.3333
$ws.Cells.Item(32, 12).Value = 1106.8462
$ws.Cells.Item(32, 13).Value = -539.3333
$ws.Cells.Item(32, 14).Value = -1758.8462
$ws.Cells.Item(33, 8).Value = 562.1142599999999
$ws.Cells.Item(33, 9).Value = 114.62963
$ws.Cells.Item(33, 11).Value = 114.62963
$ws.Cells.Item(33, 13).Value = 114.37037
$ws.Cells.Item(113, 8).Value = 1856.0869
$ws.Cells.Item(113, 9).Value = 1790
$ws.Cells.Item(113, 10).Value = 1916.6666
$ws.Cells.Item(113, 11).Value = 1790
$ws.Cells.Item(113, 12).Value = 1916.6666
$ws.Cells.Item(113, 13).Value = 1464
$ws.Cells.Item(113, 14).Value = -8424.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 917.4666999999999
$ws.Cells.Item(2, 9).Value = 773.5
$ws.Cells.Item(2, 10).Value = 1013.44446
$ws.Cells.Item(2, 11).Value = 773.5
$ws.Cells.Item(2, 12).Value = 1013.44446
$ws.Cells.Item(2, 13).Value = -660.5
$ws.Cells.Item(2, 14).Value = -1239.44446
$ws.Cells.Item(45, 8).Value = 1788.7307
$ws.Cells.Item(45, 9).Value = 1600.6957
$ws.Cells.Item(45, 10).Value = 3230.3333
$ws.Cells.Item(45, 11).Value = 1600.6957
$ws.Cells.Item(45, 12).Value = 3230.3333
$ws.Cells.Item(45, 13).Value = -1223.6957
$ws.Cells.Item(45, 14).Value = -3984.3333
$ws.Cells.Item(74, 8).Value = 58900.742
$ws.Cells.Item(74, 9).Value = 81721.52
$ws.Cells.Item(74, 10).Value = 1848.8
$ws.Cells.Item(74, 11).Value = 81721.52
$ws.Cells.Item(74, 12).Value = 1848.8
$ws.Cells.Item(74, 13).Value = -80847.52
$ws.Cells.Item(74, 14).Value = -3596.8
$ws.Cells.Item(77, 8).Value = 58900.742
$ws.Cells.Item(77, 9).Value = 81721.52
$ws.Cells.Item(77, 10).Value = 1848.8
$ws.Cells.Item(77, 11).Value = 408607.6
$ws.Cells.Item(77, 12).Value = 9244
$ws.Cells.Item(77, 13).Value = -404239.6
$ws.Cells.Item(77, 14).Value = -17980
$ws.Cells.Item(116, 8).Value = 917.4666999999999
$ws.Cells.Item(116, 9).Value = 773.5
$ws.Cells.Item(116, 10).Value = 1013.44446
$ws.Cells.Item(116, 11).Value = 773.5
$ws.Cells.Item(116, 12).Value = 1013.44446
$ws.Cells.Item(116, 13).Value = 1520.5
$ws.Cells.Item(116, 14).Value = -5601.44446
$ws.Cells.Item(122, 8).Value = 1038.7273
$ws.Cells.Item(122, 9).Value = 730.2857
$ws.Cells.Item(122, 10).Value = 1578.5
$ws.Cells.Item(122, 11).Value = 2190.8571
$ws.Cells.Item(122, 12).Value = 4735.5
$ws.Cells.Item(122, 13).Value = 259.1428999999998
$ws.Cells.Item(122, 14).Value = -9635.5
$ws.Cells.Item(132, 8).Value = 10845088
$ws.Cells.Item(132, 9).Value = 16124640
$ws.Cells.Item(132, 10).Value = 1341895.9
$ws.Cells.Item(132, 11).Value = 48373920
$ws.Cells.Item(132, 12).Value = 4025687.7
$ws.Cells.Item(132, 13).Value = -48371390
$ws.Cells.Item(132, 14).Value = -4030747.7
$ws.Cells.Item(135, 8).Value = 71223.375
$ws.Cells.Item(135, 10).Value = 71223.375
$ws.Cells.Item(135, 12).Value = 71223.375
$ws.Cells.Item(135, 14).Value = -81363.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 917.4666999999999
$ws.Cells.Item(3, 9).Value = 773.5
$ws.Cells.Item(3, 10).Value = 1013.44446
$ws.Cells.Item(3, 11).Value = 773.5
$ws.Cells.Item(3, 12).Value = 1013.44446
$ws.Cells.Item(3, 13).Value = -659.5
$ws.Cells.Item(3, 14).Value = -1241.44446
$ws.Cells.Item(22, 8).Value = 467.4
$ws.Cells.Item(22, 9).Value = 391.5
$ws.Cells.Item(22, 10).Value = 486.375
$ws.Cells.Item(22, 11).Value = 391.5
$ws.Cells.Item(22, 12).Value = 486.375
$ws.Cells.Item(22, 13).Value = -218.5
$ws.Cells.Item(22, 14).Value = -832.375
$ws.Cells.Item(134, 8).Value = 41752212
$ws.Cells.Item(134, 9).Value = 55556164
$ws.Cells.Item(134, 11).Value = 166668492
$ws.Cells.Item(134, 13).Value = -166665957

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 626.625
$ws.Cells.Item(22, 9).Value = 700.8946999999999
$ws.Cells.Item(22, 10).Value = 344.4
$ws.Cells.Item(22, 11).Value = 700.8946999999999
$ws.Cells.Item(22, 12).Value = 344.4
$ws.Cells.Item(22, 13).Value = -350.8946999999999
$ws.Cells.Item(22, 14).Value = -1044.4
$ws.Cells.Item(58, 8).Value = 849.25
$ws.Cells.Item(58, 9).Value = 708.65216
$ws.Cells.Item(58, 10).Value = 1039.4706
$ws.Cells.Item(58, 11).Value = 708.65216
$ws.Cells.Item(58, 12).Value = 1039.4706
$ws.Cells.Item(58, 13).Value = -505.65216
$ws.Cells.Item(58, 14).Value = -1445.4706
$ws.Cells.Item(105, 8).Value = 375
$ws.Cells.Item(105, 9).Value = 375
$ws.Cells.Item(105, 11).Value = 375
$ws.Cells.Item(105, 13).Value = 1372
$ws.Cells.Item(136, 8).Value = 849.25
$ws.Cells.Item(136, 9).Value = 708.65216
$ws.Cells.Item(136, 10).Value = 1039.4706
$ws.Cells.Item(136, 11).Value = 2125.95648
$ws.Cells.Item(136, 12).Value = 3118.4118
$ws.Cells.Item(136, 13).Value = 424.0435200000002
$ws.Cells.Item(136, 14).Value = -8218.4118

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 662.2857
$ws.Cells.Item(117, 9).Value = 385
$ws.Cells.Item(117, 10).Value = 1032
$ws.Cells.Item(117, 11).Value = 1155
$ws.Cells.Item(117, 12).Value = 3096
$ws.Cells.Item(117, 13).Value = 2287
$ws.Cells.Item(117, 14).Value = -9980
$ws.Cells.Item(129, 8).Value = 112341.445
$ws.Cells.Item(129, 9).Value = 553.3333
$ws.Cells.Item(129, 11).Value = 1659.9999
$ws.Cells.Item(129, 13).Value = 3340.0001
$ws.Cells.Item(131, 8).Value = 981.4464
$ws.Cells.Item(131, 9).Value = 821.5
$ws.Cells.Item(131, 10).Value = 1000.64
$ws.Cells.Item(131, 11).Value = 2464.5
$ws.Cells.Item(131, 12).Value = 3001.92
$ws.Cells.Item(131, 13).Value = 2575.5
$ws.Cells.Item(131, 14).Value = -13081.92

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 15541.893
$ws.Cells.Item(102, 9).Value = 5085.8
$ws.Cells.Item(102, 11).Value = 5085.8
$ws.Cells.Item(102, 13).Value = -3463.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 582.2778
$ws.Cells.Item(22, 9).Value = 628.7143
$ws.Cells.Item(22, 10).Value = 552.7273
$ws.Cells.Item(22, 11).Value = 628.7143
$ws.Cells.Item(22, 12).Value = 552.7273
$ws.Cells.Item(22, 13).Value = -333.7143
$ws.Cells.Item(22, 14).Value = -1142.7273
$ws.Cells.Item(27, 8).Value = 582.2778
$ws.Cells.Item(27, 9).Value = 628.7143
$ws.Cells.Item(27, 10).Value = 552.7273
$ws.Cells.Item(27, 11).Value = 628.7143
$ws.Cells.Item(27, 12).Value = 552.7273
$ws.Cells.Item(27, 13).Value = -521.7143
$ws.Cells.Item(27, 14).Value = -766.7273
$ws.Cells.Item(46, 8).Value = 1312.8572
$ws.Cells.Item(46, 9).Value = 1000
$ws.Cells.Item(46, 10).Value = 1438
$ws.Cells.Item(46, 11).Value = 1000
$ws.Cells.Item(46, 12).Value = 1438
$ws.Cells.Item(46, 13).Value = -812
$ws.Cells.Item(46, 14).Value = -1814
$ws.Cells.Item(70, 8).Value = 8720.666999999999
$ws.Cells.Item(70, 9).Value = 5000
$ws.Cells.Item(70, 10).Value = 10581
$ws.Cells.Item(70, 11).Value = 5000
$ws.Cells.Item(70, 12).Value = 10581
$ws.Cells.Item(70, 13).Value = -4730
$ws.Cells.Item(70, 14).Value = -11121
$ws.Cells.Item(73, 8).Value = 8720.666999999999
$ws.Cells.Item(73, 9).Value = 5000
$ws.Cells.Item(73, 10).Value = 10581
$ws.Cells.Item(73, 11).Value = 5000
$ws.Cells.Item(73, 12).Value = 10581
$ws.Cells.Item(73, 13).Value = -4064
$ws.Cells.Item(73, 14).Value = -12453
$ws.Cells.Item(110, 8).Value = 15828.8
$ws.Cells.Item(110, 10).Value = 15828.8
$ws.Cells.Item(110, 12).Value = 15828.8
$ws.Cells.Item(110, 14).Value = -24008.8
$ws.Cells.Item(136, 8).Value = 835888.2
$ws.Cells.Item(136, 9).Value = 2002819.6
$ws.Cells.Item(136, 10).Value = 2365.7144
$ws.Cells.Item(136, 11).Value = 6008458.800000001
$ws.Cells.Item(136, 12).Value = 7097.1432
$ws.Cells.Item(136, 13).Value = -6005908.800000001
$ws.Cells.Item(136, 14).Value = -12197.1432

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4925.3335
$ws.Cells.Item(62, 9).Value = 4825.5
$ws.Cells.Item(62, 10).Value = 5125
$ws.Cells.Item(62, 11).Value = 4825.5
$ws.Cells.Item(62, 12).Value = 5125
$ws.Cells.Item(62, 13).Value = -4201.5
$ws.Cells.Item(62, 14).Value = -6373
$ws.Cells.Item(65, 8).Value = 4925.3335
$ws.Cells.Item(65, 9).Value = 4825.5
$ws.Cells.Item(65, 10).Value = 5125
$ws.Cells.Item(65, 11).Value = 24127.5
$ws.Cells.Item(65, 12).Value = 25625
$ws.Cells.Item(65, 13).Value = -21007.5
$ws.Cells.Item(65, 14).Value = -31865
$ws.Cells.Item(132, 8).Value = 3766.7441
$ws.Cells.Item(132, 9).Value = 1180.9697
$ws.Cells.Item(132, 10).Value = 12299.8
$ws.Cells.Item(132, 11).Value = 3542.9091
$ws.Cells.Item(132, 12).Value = 36899.39999999999
$ws.Cells.Item(132, 13).Value = -1012.9091
$ws.Cells.Item(132, 14).Value = -41959.39999999999
